$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer (default, footer1.xml) Pearson logo: image1.png -> image2.png
$ftr1 = $sec.Footers.Item(1)
$shp1 = $ftr1.Range.InlineShapes.Item(1)
$shp1.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# Footer (first page, footer2.xml) Pearson logo: image1.png -> image2.png
$ftr2 = $sec.Footers.Item(2)
$shp2 = $ftr2.Range.InlineShapes.Item(1)
$shp2.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# Header (first page, header2.xml) BTec logo: image2.jpg -> image1.jpg
$hdr2 = $sec.Headers.Item(2)
$shp3 = $hdr2.Range.InlineShapes.Item(1)
$shp3.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"
